# Apply the "ch12 - ch15 airflow" update to the reading-progress tracker.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the newly-read chapters (rows 9 and 10): actual page counters
#     and tomato (pomodoro) counts for chapters 8 and 9. ---
$ws.Range("G9").Value = 301
$ws.Range("H9").Value = 312
$ws.Range("J9").Value = 3

$ws.Range("G10").Value = 313
$ws.Range("H10").Value = 393
$ws.Range("J10").Value = 11

# --- Column I ("Всего прочитано"): the day-count formula is now inclusive
#     of both the start and end day, so it gains a "+1". Re-apply it across
#     the whole used range (rows 2-10) so it covers the two new rows too. ---
$ws.Range("I2").Formula = "=H2-G2+1"
$ws.Range("I3:I8").Formula = "=H3-G3+1"
$ws.Range("I9").Formula = "=H9-G9+1"
$ws.Range("I10").Formula = "=H10-G10+1"

# --- Column K ("Стр/пом"): extend the pages-per-tomato formula down into
#     the two newly-populated rows. ---
$ws.Range("K9").Formula = "=I9/J9"
$ws.Range("K10").Formula = "=I10/J10"

# --- Restore the author's last active-cell selection. ---
$ws.Range("K14").Select()
